# Append 12 more innings rows (rows 14-25) to the Robin Uthappa match log,
# duplicating the existing 12 rows in a different order, per the commit's
# refreshed scrape of espn data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  @(14, " Abu Dhabi",   " October 30 2020",    "Royals won by 7 wickets (with 15 balls remaining)",   "Rajasthan Royals", "Kings XI Punjab",             "Robin Uthappa ", "30", "23", "1", "2", "130.43"),
  @(15, " Abu Dhabi",   " October 25 2020",    "Royals won by 8 wickets (with 10 balls remaining)",   "Rajasthan Royals", "Mumbai Indians",              "Robin Uthappa ", "13", "11", "2", "0", "118.18"),
  @(16, " Dubai (DSC)", " October 14 2020",    "Capitals won by 13 runs",                              "Rajasthan Royals", "Delhi Capitals",              "Robin Uthappa ", "32", "27", "3", "1", "118.51"),
  @(17, " Dubai (DSC)", " October 22 2020",    "Sunrisers won by 8 wickets (with 11 balls remaining)", "Rajasthan Royals", "Sunrisers Hyderabad",         "Robin Uthappa ", "19", "13", "2", "1", "146.15"),
  @(18, " Sharjah",     " September 22 2020",  "Royals won by 16 runs",                                "Rajasthan Royals", "Chennai Super Kings",         "Robin Uthappa ", "5",  "9",  "0", "0", "55.55"),
  @(19, " Dubai (DSC)", " September 30 2020",  "KKR won by 37 runs",                                   "Rajasthan Royals", "Kolkata Knight Riders",       "Robin Uthappa ", "2",  "7",  "0", "0", "28.57"),
  @(20, " Dubai (DSC)", " November 01 2020",   "KKR won by 60 runs",                                   "Rajasthan Royals", "Kolkata Knight Riders",       "Robin Uthappa ", "6",  "2",  "0", "1", "300.00"),
  @(21, " Abu Dhabi",   " October 19 2020",    "Royals won by 7 wickets (with 15 balls remaining)",   "Rajasthan Royals", "Chennai Super Kings",         "Robin Uthappa ", "4",  "9",  "0", "0", "44.44"),
  @(22, " Sharjah",     " September 27 2020",  "Royals won by 4 wickets (with 3 balls remaining)",     "Rajasthan Royals", "Kings XI Punjab",             "Robin Uthappa ", "9",  "4",  "2", "0", "225.00"),
  @(23, " Dubai (DSC)", " October 17 2020",    "RCB won by 7 wickets (with 2 balls remaining)",        "Rajasthan Royals", "Royal Challengers Bangalore", "Robin Uthappa ", "41", "22", "7", "1", "186.36"),
  @(24, " Dubai (DSC)", " October 11 2020",    "Royals won by 5 wickets (with 1 ball remaining)",      "Rajasthan Royals", "Sunrisers Hyderabad",         "Robin Uthappa ", "18", "15", "1", "1", "120.00"),
  @(25, " Abu Dhabi",   " October 03 2020",    "RCB won by 8 wickets (with 5 balls remaining)",        "Rajasthan Royals", "Royal Challengers Bangalore", "Robin Uthappa ", "17", "22", "1", "0", "77.27")
)

# Force text storage up front so numeric-looking values (e.g. "30",
# "130.43") stay strings instead of being coerced to numbers.
$fullRange = $ws.Range("A14:K25")
$fullRange.NumberFormat = "@"

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
    $ws.Range("E$r").Value = $row[5]
    $ws.Range("F$r").Value = $row[6]
    $ws.Range("G$r").Value = $row[7]
    $ws.Range("H$r").Value = $row[8]
    $ws.Range("I$r").Value = $row[9]
    $ws.Range("J$r").Value = $row[10]
    $ws.Range("K$r").Value = $row[11]
}

# Drop the temporary "@" number-format style again so the new rows don't
# pick up a distinct cell style from the original data rows.
$fullRange.ClearFormats()
